# ProjekProposal.docx edit:
#  1. Title: "Evolving Boids" -> "Evolving Clusters"
#  2. Add a new sentence mentioning the Ventrella Cluster Algorithmus right
#     after the paragraph about the awt/swing implementation libraries.

$d = $word.ActiveDocument

# --- 1. Boids -> Clusters -----------------------------------------------
$d.Content.Find.Execute("Boids", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Clusters", 2) | Out-Null

# --- 2. Insert the new paragraph about the Ventrella Cluster Algorithmus -
$inserted = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*vorallem mit der awt und auch swing Library implementiert werden.*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "Als Algorithmus für die Partikelsimulation wird der Ventrella Cluster Algorithmus verwendet."
        $inserted = $true
        break
    }
}
if (-not $inserted) {
    Write-Output "WARNING: target paragraph for Ventrella insertion not found"
}
